$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline / page-setup / margin conventions used by the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 3. Copy header-row formatting (bold/centered/bordered) and date-column
#        formatting from the "Weekly Quantity" sheet so the new sheet matches
#        the workbook's existing look-and-feel ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A76").PasteSpecial(-4122)

# --- 4. Header row values ---
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# --- 5. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @"
45011.99999999999,6,-12.0371285290029,21.30696095054402
45025.99999999999,6,-10.73247897206887,22.38573317689519
45032.99999999999,6,-10.99924221458364,22.07740565974546
45039.99999999999,6,-12.37545363725458,22.67110176692979
45046.99999999999,6,-9.794621070320954,23.24612333843451
45053.99999999999,6,-10.14327721574721,23.26786519345735
45060.99999999999,6,-10.6046371734621,23.9513044300436
45067.99999999999,7,-10.068416199098,21.9584942327692
45074.99999999999,7,-9.670733561758894,23.49639746022065
45088.99999999999,7,-9.488562021826844,25.26087377334528
45095.99999999999,7,-8.545993440543349,23.84767213673662
45102.99999999999,7,-10.00163483418195,23.64620113616784
45116.99999999999,8,-7.936356750355466,24.70754042263821
45123.99999999999,8,-10.1044238940701,24.35052551170137
45130.99999999999,8,-9.751837846937343,24.61683008328276
45137.99999999999,8,-8.431478319955771,24.39360937612444
45144.99999999999,8,-10.23235838350582,24.38701480939397
45151.99999999999,8,-8.674592886312524,25.38536453678897
45158.99999999999,8,-8.972765552059865,25.40885950522467
45165.99999999999,9,-9.178825777156257,26.6373947820951
45172.99999999999,9,-6.372171696929078,26.29239561238792
45179.99999999999,9,-7.814319573920327,25.62695833523329
45186.99999999999,9,-6.526330067917297,27.08137352270162
45193.99999999999,9,-7.023502018803199,25.70756091082966
45207.99999999999,9,-7.210156528131384,26.02414732228687
45214.99999999999,10,-6.939774794720397,26.34029451348136
45221.99999999999,10,-7.754918779823782,26.73777274780085
45228.99999999999,10,-5.923564411138019,26.69647454460682
45235.99999999999,10,-7.348934237936498,26.06398989750172
45242.99999999999,10,-7.172799404431317,27.0872270278123
45249.99999999999,10,-5.617534263584985,26.35434301199912
45256.99999999999,10,-5.591987872018462,27.13450319817693
45263.99999999999,11,-6.267188011361647,27.33451519834306
45270.99999999999,11,-6.161775081884848,27.58956231966536
45277.99999999999,11,-6.186207923318476,26.94589698379255
45298.99999999999,11,-5.284999898842883,28.07611161408637
45312.99999999999,12,-4.95500712593515,28.49682547825615
45319.99999999999,12,-4.975133677987678,28.13517347342118
45326.99999999999,12,-4.793826456843298,29.07524094114235
45333.99999999999,12,-4.612586857643346,28.28538307242122
45340.99999999999,12,-3.455670603357186,27.75975144474019
45347.99999999999,12,-4.343218925539384,30.16407085822492
45361.99999999999,13,-4.276504261486314,30.50523426904406
45368.99999999999,13,-4.083368895071146,29.40488225272916
45375.99999999999,13,-3.61387940221831,28.9523804988183
45382.99999999999,13,-2.840239776790965,29.13735000058245
45389.99999999999,13,-4.334664772405904,30.26614595473733
45396.99999999999,13,-4.005522800327144,30.45970974999945
45403.99999999999,13,-2.96430987988622,30.03560724846849
45410.99999999999,14,-3.538399503552816,31.33636499493482
45417.99999999999,14,-4.374326201976882,30.39045569726521
45424.99999999999,14,-2.741423024720475,31.4028842268829
45431.99999999999,14,-1.608538533234486,30.96295800967486
45438.99999999999,14,-2.316780791531706,30.38706353251931
45445.99999999999,14,-1.861715848723664,30.85957003701709
45452.99999999999,14,-1.563022847385252,31.47627028039931
45459.99999999999,15,-2.367213101294397,32.08318566034423
45473.99999999999,15,-2.199582385742347,31.99453918391351
45480.99999999999,15,-2.192669765037916,32.52577640080861
45487.99999999999,15,-1.816070903057942,31.46862238222596
45494.99999999999,15,-1.260638802646808,32.25897164193994
45501.99999999999,15,-0.9438240231370745,32.19249826253025
45508.99999999999,16,-1.809895504640155,32.13950177895764
45515.99999999999,16,-1.265779355088546,32.94492601383711
45529.99999999999,16,-0.4977069914080139,31.63373838064538
45536.99999999999,16,0.4178791228183382,33.07292975631397
45543.99999999999,16,-0.2998805578585826,34.09479120450891
45550.99999999999,16,0.983888692075234,31.64861286103812
45557.99999999999,17,0.5870641757021132,33.85986017056776
45564.99999999999,17,-1.028476336887851,34.21179423218596
45571.99999999999,17,0.1267984850284437,33.77498989525646
45578.99999999999,17,-1.933593172804086,32.69294326547109
45585.99999999999,17,0.399932403916532,34.08908814309866
45592.99999999999,17,0.9002828337590718,33.07792581988232
45599.99999999999,17,1.842973943069506,34.26760263977399
"@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $wsForecast.Cells.Item($r, 1).Value = [double]$parts[0]
    $wsForecast.Cells.Item($r, 2).Value = [double]$parts[1]
    $wsForecast.Cells.Item($r, 3).Value = [double]$parts[2]
    $wsForecast.Cells.Item($r, 4).Value = [double]$parts[3]
    $r = $r + 1
}

Write-Output "Rows written: $($r - 2)"
